# Split the "15 mins" run into its own run and correct it to "15 min" in the
# "each group will present the work (15 mins), preparing a set of slides..."
# sentence found on the "Group Project (HOW ?)" slide.
#
# Before:
#   <a:r><a:t>each group will present the work (15 mins), preparing a set of
#             slides to comment the results</a:t></a:r>
# After:
#   <a:r><a:t>each group will present the work (</a:t></a:r>
#   <a:r><a:t>15 min), </a:t></a:r>
#   <a:r><a:t>preparing a set of slides to comment the results</a:t></a:r>

$p = $ppt.ActivePresentation

$oldFragment = "15 mins), "
$newFragment = "15 min), "
$searchHint  = "each group will present the work"

# Locate the shape holding the sentence without hard-coding slide/shape
# indices, so the edit keeps working even if the deck is reordered.
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                $text = $shape.TextFrame.TextRange.Text
                if ($text -like "*$searchHint*") {
                    $targetShape = $shape
                }
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not find the shape containing '$searchHint'"
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text
$startPos = $fullText.IndexOf($oldFragment)
if ($startPos -lt 0) {
    throw "Could not find '$oldFragment' in the target shape's text"
}

# Characters() is 1-indexed; grab just the "15 mins), " span so it becomes
# its own run (the surrounding text keeps its existing runs untouched).
$middle = $tr.Characters($startPos + 1, $oldFragment.Length)
$middle.Text = $newFragment
